# Update gh-pages to output generated at 14f1f32
# Update attendance counts (column F) and add bilibili detail links (column I)
# for both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (first worksheet) ---
$ws1.Range("F2").Value = 2076
$ws1.Range("I2").Value = "https://show.bilibili.com/platform/detail.html?id=79764"

$ws1.Range("F3").Value = 609
$ws1.Range("I3").Value = "https://show.bilibili.com/platform/detail.html?id=78089"

$ws1.Range("F4").Value = 1428
$ws1.Range("I4").Value = "https://show.bilibili.com/platform/detail.html?id=79354"

$ws1.Range("F5").Value = 6940
$ws1.Range("I5").Value = "https://show.bilibili.com/platform/detail.html?id=77938"

$ws1.Range("F6").Value = 172
$ws1.Range("I6").Value = "https://show.bilibili.com/platform/detail.html?id=79051"

$ws1.Range("F7").Value = 101
$ws1.Range("I7").Value = "https://show.bilibili.com/platform/detail.html?id=80943"

# --- Sheet "全部类型" (fourth worksheet) ---
$ws4.Range("F2").Value = 2076
$ws4.Range("I2").Value = "https://show.bilibili.com/platform/detail.html?id=79764"

$ws4.Range("F3").Value = 609
$ws4.Range("I3").Value = "https://show.bilibili.com/platform/detail.html?id=78089"

$ws4.Range("F4").Value = 1429
$ws4.Range("I4").Value = "https://show.bilibili.com/platform/detail.html?id=79354"

$ws4.Range("F5").Value = 6940
$ws4.Range("I5").Value = "https://show.bilibili.com/platform/detail.html?id=77938"

$ws4.Range("F6").Value = 172
$ws4.Range("I6").Value = "https://show.bilibili.com/platform/detail.html?id=79051"

$ws4.Range("F7").Value = 101
$ws4.Range("I7").Value = "https://show.bilibili.com/platform/detail.html?id=80943"
